$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cols = @("B","C","E","F","G","H","N","O")

$data = @{
    2 = @(0.7270418549378519, 0.2183220345354471, 0.4288618641543138, 0.4443680307746121, 0.1658682458822369, 0.3179217265166372, 0.7759822506049687, 0.8688947377034424)
    3 = @(0.6344550015893731, 0.1986302810215079, 0.3740835637659217, 0.3878228170618172, 0.161509437225817, 0.3197718048556268, 0.7715108989491313, 0.8632572353802885)
    4 = @(0.5774251797177499, 0.1864711635582523, 0.3405328773521035, 0.3531389305169483, 0.1590513841166938, 0.3211357469931855, 0.769101914212726, 0.8607265276269374)
    5 = @(0.5541407890677021, 0.1814993588170921, 0.3268800212226495, 0.3390132514313251, 0.1581041564951278, 0.3217487728211879, 0.768205152982361, 0.85992809931804)
    6 = @(0.5502717971138793, 0.1806727849171068, 0.3246140936602728, 0.336668177824194, 0.1579501458493269, 0.3218540170123916, 0.7680613848783224, 0.8598095501004792)
    7 = @(0.5771113358654532, 0.1864041799372842, 0.3403486740989905, 0.3529483938344953, 0.1590383895920198, 0.3211437830013395, 0.7690894759167719, 0.8607148184967173)
    8 = @(0.6951562310368331, 0.2115466731633262, 0.4099559814636962, 0.4248636149813478, 0.1643197610470892, 0.3185122416545596, 0.7743708854088425, 0.8667570323315346)
    9 = @(0.9251601606953841, 0.2602977884466497, 0.5471985066635767, 0.5661985755041457, 0.1764288983113573, 0.3151664844223063, 0.7873849075876933, 0.8860465031119986)
    10 = @(1.09319709068825, 0.2957661528352844, 0.6486047476104773, 0.6702781546542269, 0.186423211927746, 0.3138229417580192, 0.7985520518447515, 0.9048347211059422)
    11 = @(1.169427645025678, 0.3118235522812256, 0.6948885853739597, 0.7176906081379002, 0.1912146412897897, 0.3134555769324265, 0.8039782092016026, 0.9144018107482168)
    12 = @(1.198262884874339, 0.3178926856266173, 0.7124391560133745, 0.7356546913071611, 0.1930647511803016, 0.313351670539106, 0.8060824742386217, 0.9181726839399573)
    13 = @(1.192054128057464, 0.3165861039501863, 0.7086582355407529, 0.7317853510981394, 0.1926647016155272, 0.3133724804842188, 0.8056270861294905, 0.917353954549867)
    14 = @(1.171800579071487, 0.3123230950230891, 0.6963319923408449, 0.7191683204515869, 0.1913661327561869, 0.3134463220586667, 0.8041503379472914, 0.9147090685915771)
    15 = @(1.159390537327965, 0.3097103772528556, 0.6887849766281704, 0.7114413442032514, 0.1905753850408018, 0.3134961414438493, 0.8032522247780776, 0.9131083157428179)
    16 = @(1.088210895231327, 0.2947151791693159, 0.6455832343717702, 0.6671810134426437, 0.1861150505825151, 0.3138518686272391, 0.7982043812617121, 0.9042301314289034)
    17 = @(1.044489648070794, 0.2854960643796858, 0.6191209843935326, 0.6400460337125793, 0.1834418585014816, 0.3141326472300818, 0.7951961521616084, 0.8990458948276512)
    18 = @(1.019322631669809, 0.2801862138899764, 0.6039149503132393, 0.6244449056556647, 0.181927338129185, 0.3143170843340641, 0.793498497650333, 0.8961600096104121)
    19 = @(1.010798168984877, 0.2783871531427167, 0.5987688590430622, 0.619163680173358, 0.1814184883805297, 0.3143834672253973, 0.7929293086531572, 0.8951993380266003)
    20 = @(1.04914589952881, 0.2864782086466846, 0.6219364377537886, 0.642933953830422, 0.183724037514466, 0.3141003825390953, 0.7955130114367108, 0.8995878256149581)
    21 = @(1.177750406370478, 0.313575557631026, 0.6999518461297498, 0.7228739723491628, 0.191746581253625, 0.3134236763335281, 0.8045827535602541, 0.9154819078869707)
    22 = @(1.261615965852798, 0.3312183174540451, 0.7510794944355155, 0.7751780083420101, 0.1971980819835295, 0.3131866894222526, 0.8107987349008141, 0.9267330186176821)
    23 = @(1.216872711842257, 0.3218082780102804, 0.7237782928742718, 0.7472568307830727, 0.1942692969134612, 0.3132943423467793, 0.8074548531010635, 0.9206486518461361)
    24 = @(1.047040904981429, 0.286034211287074, 0.6206635482033107, 0.6416283278902171, 0.18359639485125, 0.3141148977271087, 0.7953696602267684, 0.8993425241202715)
    25 = @(0.8631011742364763, 0.2471697761009466, 0.5099787580083301, 0.5279251897347166, 0.1729624105080205, 0.3158764701235413, 0.7835812886509785, 0.8800233311974921)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
